{"js": "// Remove the trailing \"Ver no Jupiter...\" / copyright footer block that\n// followed the last \"Requisitos\" entry (LOQ4083), along with the blank\n// paragraph that separated it from that entry.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nfor (let i = items.length - 1; i >= 0; i--) {\n  const text = items[i].text;\n  if (text === \"Ver no Jupiter Salvar em pdf Salvar em docx\") {\n    // Delete the copyright paragraph right after this one (if present).\n    if (i + 1 < items.length) {\n      items[i + 1].delete();\n    }\n    // Delete the \"Ver no Jupiter...\" paragraph itself.\n    items[i].delete();\n    // Delete the blank separator paragraph right before this one.\n    if (i - 1 >= 0) {\n      items[i - 1].delete();\n    }\n    break;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter...\" / copyright footer block that\n# followed the last \"Requisitos\" entry (LOQ4083), along with the blank\n# paragraph that separated it from that entry.\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\nfor ($i = $count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    $text = $p.Range.Text\n\n    if ($text.TrimEnd(\"`r`a`n\") -eq \"Ver no Jupiter Salvar em pdf Salvar em docx\") {\n        # Delete the copyright paragraph that immediately follows, if present.\n        if (($i + 1) -le $d.Paragraphs.Count) {\n            $d.Paragraphs.Item($i + 1).Range.Delete()\n        }\n\n        # Delete the \"Ver no Jupiter...\" paragraph itself.\n        $d.Paragraphs.Item($i).Range.Delete()\n\n        # Delete the blank separator paragraph right before it.\n        if (($i - 1) -ge 1) {\n            $d.Paragraphs.Item($i - 1).Range.Delete()\n        }\n\n        break\n    }\n}\n"}
